$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E8"  = 16.11280000000001
    "E10" = 16.47669999999998
    "E12" = 18.11140000000002
    "E18" = 17.6355
    "E37" = 16.5796
    "E55" = 16.6076
    "E68" = 17.09180000000001
    "E77" = 18.29680000000002
    "E78" = 16.62860000000003
    "E81" = 16.5241
    "E82" = 16.6022
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
